$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.243.66"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.860.97"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.72"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4682"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2866"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06543"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.70"
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07927"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.62"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "1.865.82"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.178"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6806"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.71"
$ws.Range("E16").Value = "  -5.34%  "
$ws.Range("D17").Value = "30.238.72"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.73"
$ws.Range("E18").Value = "  +8.08%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007405"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "2.112.76"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.192"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.08"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.236"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.969"
$ws.Range("E28").Value = "  +2.75%  "
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09885"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.393"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.071"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04695"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7041"
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01881"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.637"
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.247"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.40"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8481"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4171"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9998"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.47"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "962.67"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.157"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.229"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.18"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05660"
$ws.Range("E51").Value = "  +0.40%  "
